$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the three "nophoto" headshot URL cells (rows for Donovan Clingan,
# Yarin Hasson, Alex Karaban) leaving the cell empty but keeping its
# existing text-format style.
$ws.Range("C17").ClearContents() | Out-Null
$ws.Range("C19").ClearContents() | Out-Null
$ws.Range("C26").ClearContents() | Out-Null

# Update the view: scroll back to the top-left and select a single cell.
$ws.Range("C12").Select() | Out-Null
